{"js": "// The document ends with a paragraph that only carries the \"_GoBack\"\n// bookmark, followed by a final empty, centered paragraph. The edit:\n//   1) inserts two empty, centered paragraphs right before the bookmark\n//      paragraph,\n//   2) adds the text \"T4s\" (and centers the paragraph) on the bookmark\n//      paragraph itself, keeping the bookmark tags in place,\n//   3) inserts one more empty paragraph right after the bookmark\n//      paragraph (before the pre-existing trailing centered paragraph).\n\nconst body = context.document.body;\n\n// Locate the paragraph that hosts the \"_GoBack\" bookmark robustly\n// (rather than assuming a fixed index).\nconst bookmarkRange = body.getBookmarkRange(\"_GoBack\");\nconst bookmarkPara = bookmarkRange.paragraphs.getFirst();\nbookmarkPara.load(\"text\");\nawait context.sync();\n\n// 1) Insert the trailing empty paragraph first, while the bookmark\n//    paragraph still has its original (non-centered) formatting, so the\n//    newly inserted paragraph does not inherit centered alignment.\nconst afterPara = bookmarkPara.insertParagraph(\"\", \"After\");\n\n// 2) Insert two empty paragraphs before the bookmark paragraph and\n//    center them.\nconst blank1 = bookmarkPara.insertParagraph(\"\", \"Before\");\nblank1.alignment = Word.Alignment.centered;\n\nconst blank2 = bookmarkPara.insertParagraph(\"\", \"Before\");\nblank2.alignment = Word.Alignment.centered;\n\n// 3) Add the \"T4s\" text at the start of the bookmark paragraph (keeping\n//    the bookmark itself intact) and center that paragraph too.\nbookmarkPara.insertText(\"T4s\", \"Start\");\nbookmarkPara.alignment = Word.Alignment.centered;\n\nawait context.sync();\n", "ps1": "# The document ends with a paragraph that only carries the \"_GoBack\"\n# bookmark, followed by a final empty, centered paragraph. The edit:\n#   1) inserts two empty, centered paragraphs right before the bookmark\n#      paragraph,\n#   2) adds the text \"T4s\" (and centers the paragraph) on the bookmark\n#      paragraph itself, keeping the bookmark tags in place,\n#   3) inserts one more empty paragraph right after the bookmark\n#      paragraph (before the pre-existing trailing centered paragraph).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that hosts the \"_GoBack\" bookmark robustly\n# (rather than assuming a fixed index).\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bookmarkPara = $bm.Range.Paragraphs.Item(1)\n\n# 1) Insert the trailing empty paragraph first, while the bookmark\n#    paragraph still has its original (non-centered) formatting, so the\n#    newly inserted paragraph does not inherit centered alignment.\n$bookmarkPara.Range.InsertParagraphAfter()\n\n# 2) Insert two empty paragraphs before the bookmark paragraph.\n$bookmarkPara.Range.InsertParagraphBefore()\n$bookmarkPara.Range.InsertParagraphBefore()\n\n# Paragraph indices shifted because of the inserts above, so re-resolve\n# the bookmark paragraph through the bookmark again instead of reusing a\n# stale reference/index.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bookmarkPara = $bm.Range.Paragraphs.Item(1)\n$bmIdx = $bookmarkPara.Index\n\n$d.Paragraphs.Item($bmIdx - 2).Alignment = 1\n$d.Paragraphs.Item($bmIdx - 1).Alignment = 1\n\n# 3) Add \"T4s\" text at the start of the bookmark paragraph (keeping the\n#    bookmark itself intact) and center that paragraph too.\n$bookmarkPara.Range.InsertBefore(\"T4s\")\n$bookmarkPara.Alignment = 1\n"}
